$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# personnel sheet: shrink "username" length, change "password" datatype to
# "binary" and its length.
# ---------------------------------------------------------------------------
$wsPersonnel = $wb.Worksheets.Item("personnel")
$wsPersonnel.Range("C3").Value = 50
$wsPersonnel.Range("B4").Value = "binary"
$wsPersonnel.Range("C4").Value = 64

# ---------------------------------------------------------------------------
# location sheet: a handful of length / index tweaks.
# ---------------------------------------------------------------------------
$wsLocation = $wb.Worksheets.Item("location")
$wsLocation.Range("C2").Value = 2
$wsLocation.Range("E5").Value = $false
$wsLocation.Range("C6").Value = 5
$wsLocation.Range("C8").Value = 40

# ---------------------------------------------------------------------------
# palette sheet: widen the "foreign key" column, flip an index flag, and
# bump a length.
# ---------------------------------------------------------------------------
$wsPalette = $wb.Worksheets.Item("palette")
$wsPalette.Columns.Item(7).ColumnWidth = 57.5
$wsPalette.Range("E9").Value = $true
$wsPalette.Range("C14").Value = 40

# ---------------------------------------------------------------------------
# Selections / active sheet. Replay them in the same order the author would
# have clicked through the tabs so the final active tab & per-sheet cursor
# positions all land correctly.
# ---------------------------------------------------------------------------
$wsType = $wb.Worksheets.Item("type")
$wsType.Activate()
$wsType.Range("A1:I1").Select()

$wsPermissions = $wb.Worksheets.Item("permissions")
$wsPermissions.Activate()
$wsPermissions.Range("I3").Select()

$wsLocation.Activate()
$wsLocation.Range("E2").Select()

$wsPalette.Activate()
$wsPalette.Range("E10").Select()

$wsPersonnel.Activate()
$wsPersonnel.Range("I12").Select()
